$wb = $excel.ActiveWorkbook

# "General" sheet: merge "Translation Owner ID" / "Translation Owner Name"
# into a single "Translation Owner" row.
$general = $wb.Worksheets.Item("General")
$general.Rows.Item(3).Delete()
$general.Cells.Item(2, 1).Value = "Translation Owner"
$general.Range("B8").Select()

# "Attributes" sheet: move the active selection to B2.
$attributes = $wb.Worksheets.Item("Attributes")
$attributes.Range("B2").Select()
